# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C15").Value = 45243
